# The workbook gains a new "maskRamp" column that is inserted right before
# the existing "cont1_maskOn" column (current column O), pushing the three
# trailing columns (cont1_maskOn, cont0_maskOn, cont1_maskOff) one place to
# the right (O->P, P->Q, Q->R).  A few data values are also updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column O; this shifts O:Q -> P:R and extends
# the used range out to column R.
$ws.Columns("O").Insert()

# New column O: "maskRamp" header + data (all zero for the existing rows).
$ws.Range("O1").Value = "maskRamp"
$ws.Range("O2:O5").Value = 0

# The old "cont1_maskOn" column (now column P) gets new data values.
$ws.Range("P2:P5").Value = 12

# stimT (column F) changes from 1000 to 2000 for every condition row.
$ws.Range("F2:F5").Value = 2000

# Update the active selection to match the saved view state.
$null = $ws.Range("P6").Select()
